$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-11 Sunday" "2025-05-12 Monday"

Replace-Text "148×3=444" "483×7=3381"
Replace-Text "715×2=1430" "619×8=4952"
Replace-Text "266×7=1862" "776×6=4656"
Replace-Text "671×7=4697" "959×7=6713"
Replace-Text "599×9=5391" "262×7=1834"

Replace-Text "402×5=2010" "698×5=3490"
Replace-Text "593×2=1186" "673×6=4038"
Replace-Text "166×5=830" "386×3=1158"
Replace-Text "207×7=1449" "598×4=2392"
Replace-Text "343×8=2744" "824×9=7416"

Replace-Text "380×2=760" "552×9=4968"
Replace-Text "488×9=4392" "568×2=1136"
Replace-Text "521×6=3126" "985×9=8865"
Replace-Text "509×4=2036" "498×3=1494"
Replace-Text "227×6=1362" "273×9=2457"

Replace-Text "928×2=1856" "118×4=472"
Replace-Text "785×7=5495" "623×8=4984"
Replace-Text "196×6=1176" "850×2=1700"
Replace-Text "934×7=6538" "331×7=2317"
Replace-Text "678×7=4746" "441×2=882"

Replace-Text "978×2=1956" "418×6=2508"
Replace-Text "950×4=3800" "838×7=5866"
Replace-Text "901×5=4505" "963×4=3852"
Replace-Text "851×4=3404" "561×3=1683"
Replace-Text "996×9=8964" "264×7=1848"
